# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook is a single "Estado de Cuenta" (account statement) report for
# NIT 9000855752. This edit:
#   1. Updates the aggregate "VALOR MORA" total (E11).
#   2. Updates the worker/period counters (C13, F13).
#   3. Replaces the first worker's data (rows 16-17, previously BEATRIZ BRAVO
#      PACHECO / doc 45526840) with the other worker already present further
#      down (MARIA TERESA MARMOL BARBOZA / doc 1143346806), and updates the
#      "Salario Basico" amounts.
#   4. Adds a third period row (2508) for that worker, re-using the visual
#      styling of the (now redundant) last data row, then removes that
#      now-duplicate row - shifting the signature block rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Aggregate "VALOR MORA" total ---------------------------------------
$ws.Range("E11").Value = 170820

# --- 2. Worker / period counters -------------------------------------------
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3

# --- 3. Replace data rows 16 & 17 with MARIA TERESA MARMOL BARBOZA's data --
$ws.Range("C16").Value = "1143346806"
$ws.Range("D16").Value = "MARIA TERESA MARMOL BARBOZA"
$ws.Range("E16").Value = "2506"
$ws.Range("G16").Value = 1423500

$ws.Range("C17").Value = "1143346806"
$ws.Range("D17").Value = "MARIA TERESA MARMOL BARBOZA"
$ws.Range("E17").Value = "2507"
$ws.Range("G17").Value = 1423500

# --- 4. Add the new "2508" period row ---------------------------------------
# Row 18 currently duplicates row 17's content with the lighter (top) border
# styling; row 19 already holds the same worker's data with the closing
# (bottom) border styling that the new last row needs. Copy that formatting
# onto row 18, fill in the new period, then drop the now-superfluous row 19
# (this also shifts the signature block, old rows 24-25, up to 23-24 to close
# the gap - matching the target layout).
$ws.Range("B19:J19").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

$ws.Range("C18").Value = "1143346806"
$ws.Range("D18").Value = "MARIA TERESA MARMOL BARBOZA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Rows.Item(19).Delete()
